$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("J2").Value = 3.95
$ws.Range("P2").Value = 1.9
$ws.Range("R2").Value = 1.35
$ws.Range("T2").Value = 1.98

# Row 3
$ws.Range("F3").Value = 1.23
$ws.Range("G3").Value = 1.26
$ws.Range("N3").Value = 8.6
$ws.Range("T3").Value = 1.78
$ws.Range("AG3").Value = 13.5

# Row 4
$ws.Range("I4").Value = 2.26
$ws.Range("J4").Value = 3.35
$ws.Range("K4").Value = 3.4
$ws.Range("M4").Value = 1.11
$ws.Range("N4").Value = 3
$ws.Range("O4").Value = 1.48
$ws.Range("P4").Value = 1.67
$ws.Range("T4").Value = 2.06
$ws.Range("U4").Value = 1.89
$ws.Range("V4").Value = 1.8
$ws.Range("W4").Value = 1.34
$ws.Range("AC4").Value = 7.4

# Row 5
$ws.Range("F5").Value = 2.76
$ws.Range("G5").Value = 2.8
$ws.Range("K5").Value = 3.35
$ws.Range("L5").Value = 1.46
$ws.Range("P5").Value = 1.71
$ws.Range("Q5").Value = 2.32
$ws.Range("V5").Value = 1.51
$ws.Range("W5").Value = 1.55
$ws.Range("X5").Value = 10.5
$ws.Range("Z5").Value = 17.5
$ws.Range("AA5").Value = 50
$ws.Range("AC5").Value = 7.2
$ws.Range("AD5").Value = 13
$ws.Range("AE5").Value = 36
$ws.Range("AF5").Value = 16.5
$ws.Range("AG5").Value = 12.5
$ws.Range("AH5").Value = 19.5
$ws.Range("AI5").Value = 60
$ws.Range("AJ5").Value = 42
$ws.Range("AK5").Value = 34
$ws.Range("AL5").Value = 55
$ws.Range("AM5").Value = 140
$ws.Range("AN5").Value = 34
$ws.Range("AO5").Value = 40

# Row 6
$ws.Range("L6").Value = 1.44
$ws.Range("P6").Value = 1.82
$ws.Range("V6").Value = 1.62
$ws.Range("W6").Value = 1.44
$ws.Range("X6").Value = 11.5
$ws.Range("Y6").Value = 10
$ws.Range("Z6").Value = 15.5
$ws.Range("AA6").Value = 36
$ws.Range("AB6").Value = 11.5
$ws.Range("AD6").Value = 11.5
$ws.Range("AE6").Value = 29
$ws.Range("AF6").Value = 20
$ws.Range("AG6").Value = 13.5
$ws.Range("AH6").Value = 18.5
$ws.Range("AI6").Value = 44
$ws.Range("AK6").Value = 38
$ws.Range("AM6").Value = 110
$ws.Range("AN6").Value = 38
$ws.Range("AO6").Value = 26
